$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row keeps the same two labels; just re-assert them so shared-string
# order is driven by first-use order below (matches target sst ordering).
$ws.Range("A1").Value = "Cód. tema"
$ws.Range("B1").Value = "Texto pregunta"

# Rows 2-5: replace topic codes + sample questions with the new lorem-ipsum set
$ws.Range("A2").Value = "en7-u01"
$ws.Range("B2").Value = "¿Lorem ipsum dolor sit amet, consectetur adipiscing elit. In mi ante, ultricies eu massa vitae, maximus faucibus nisl.?"

$ws.Range("A3").Value = "en7-u01"
$ws.Range("B3").Value = "¿Nulla efficitur purus a risus gravida, quis mattis urna semper?"

$ws.Range("B4").Value = "¿Aenean bibendum ante vitae mi ultrices, nec dictum nibh vulputate. Maecenas fermentum blandit nibh, a cursus risus lobortis eu?"

$ws.Range("B5").Value = "¿In gravida cursus quam at aliquam. Sed erat nunc, accumsan a gravida ac, maximus ut est?"

$ws.Range("A4").Value = "en7-u02"
$ws.Range("A5").Value = "en7-u02"

# Row 6 ("l1-u01" in A6 only) is wiped out completely - no cells remain in it
$ws.Range("A6").ClearContents()

# Row 7 keeps only the styled, now-empty B7 cell - drop A7 and B7's text/value
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()

# Move/save the active selection to A6 (it used to be B7)
$ws.Range("A6").Select()
